$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 <- original row 17 data (columns B:AC)
$ws.Range("B16").Value = 6799836
$ws.Range("F16").Value = "Maccabi Tel Aviv"
$ws.Range("G16").Value = "Maccabi Bnei Raina"
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = "D"
$ws.Range("K16").Value = 1.181
$ws.Range("L16").Value = 6
$ws.Range("M16").Value = 11
$ws.Range("N16").Value = 1.2
$ws.Range("O16").Value = 6
$ws.Range("P16").Value = 10
$ws.Range("Q16").Value = -1.75
$ws.Range("R16").Value = 1.825
$ws.Range("S16").Value = 2.025
$ws.Range("T16").Value = 3
$ws.Range("U16").Value = 1.85
$ws.Range("V16").Value = 2
$ws.Range("W16").Value = -1
$ws.Range("X16").Value = 5
$ws.Range("Y16").Value = -1
$ws.Range("Z16").Value = -1
$ws.Range("AA16").Value = 1.025
$ws.Range("AB16").Value = -1
$ws.Range("AC16").Value = 1

# Row 17 <- original row 16 data (columns B:AC)
$ws.Range("B17").Value = 6799838
$ws.Range("F17").Value = "Hapoel Haifa"
$ws.Range("G17").Value = "MS Ashdod"
$ws.Range("H17").Value = 2
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = "H"
$ws.Range("K17").Value = 2.15
$ws.Range("L17").Value = 3
$ws.Range("M17").Value = 3.2
$ws.Range("N17").Value = 2.15
$ws.Range("O17").Value = 3.1
$ws.Range("P17").Value = 3.1
$ws.Range("Q17").Value = -0.25
$ws.Range("R17").Value = 2
$ws.Range("S17").Value = 1.85
$ws.Range("T17").Value = 2.5
$ws.Range("U17").Value = 2
$ws.Range("V17").Value = 1.85
$ws.Range("W17").Value = 1.15
$ws.Range("X17").Value = -1
$ws.Range("Y17").Value = -1
$ws.Range("Z17").Value = 1
$ws.Range("AA17").Value = -1
$ws.Range("AB17").Value = -1
$ws.Range("AC17").Value = 0.8500000000000001

# Row 18 <- original row 19 data (columns B:AC)
$ws.Range("B18").Value = 6799832
$ws.Range("F18").Value = "Maccabi Petach Tikva"
$ws.Range("G18").Value = "Maccabi Haifa"
$ws.Range("H18").Value = 3
$ws.Range("I18").Value = 2
$ws.Range("J18").Value = "H"
$ws.Range("K18").Value = 5
$ws.Range("L18").Value = 3.75
$ws.Range("M18").Value = 1.533
$ws.Range("N18").Value = 4.75
$ws.Range("O18").Value = 3.6
$ws.Range("P18").Value = 1.571
$ws.Range("Q18").Value = 0.75
$ws.Range("R18").Value = 2.05
$ws.Range("S18").Value = 1.8
$ws.Range("T18").Value = 2.5
$ws.Range("U18").Value = 1.825
$ws.Range("V18").Value = 2.025
$ws.Range("W18").Value = 3.75
$ws.Range("X18").Value = -1
$ws.Range("Y18").Value = -1
$ws.Range("Z18").Value = 1.05
$ws.Range("AA18").Value = -1
$ws.Range("AB18").Value = 0.825
$ws.Range("AC18").Value = -1

# Row 19 <- original row 18 data (columns B:AC)
$ws.Range("B19").Value = 6799834
$ws.Range("F19").Value = "Maccabi Netanya"
$ws.Range("G19").Value = "Hapoel Petah Tikva"
$ws.Range("H19").Value = 4
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = "H"
$ws.Range("K19").Value = 1.571
$ws.Range("L19").Value = 3.6
$ws.Range("M19").Value = 5
$ws.Range("N19").Value = 1.65
$ws.Range("O19").Value = 3.6
$ws.Range("P19").Value = 4.75
$ws.Range("Q19").Value = -0.75
$ws.Range("R19").Value = 1.925
$ws.Range("S19").Value = 1.925
$ws.Range("T19").Value = 2.5
$ws.Range("U19").Value = 1.875
$ws.Range("V19").Value = 1.975
$ws.Range("W19").Value = 0.6499999999999999
$ws.Range("X19").Value = -1
$ws.Range("Y19").Value = -1
$ws.Range("Z19").Value = 0.925
$ws.Range("AA19").Value = -1
$ws.Range("AB19").Value = 0.875
$ws.Range("AC19").Value = -1

# Row 58 <- original row 60 data (columns B:AC)
$ws.Range("B58").Value = 7542748
$ws.Range("F58").Value = "MS Ashdod"
$ws.Range("G58").Value = "Hapoel Jerusalem FC"
$ws.Range("H58").Value = 2
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = "H"
$ws.Range("K58").Value = 2.5
$ws.Range("L58").Value = 3.2
$ws.Range("M58").Value = 2.625
$ws.Range("N58").Value = 2.4
$ws.Range("O58").Value = 2.9
$ws.Range("P58").Value = 3
$ws.Range("Q58").Value = -0.25
$ws.Range("R58").Value = 2.125
$ws.Range("S58").Value = 1.75
$ws.Range("T58").Value = 2
$ws.Range("U58").Value = 2.05
$ws.Range("V58").Value = 1.8
$ws.Range("W58").Value = 1.4
$ws.Range("X58").Value = -1
$ws.Range("Y58").Value = -1
$ws.Range("Z58").Value = 1.125
$ws.Range("AA58").Value = -1
$ws.Range("AB58").Value = 0
$ws.Range("AC58").Value = 0

# Row 60 <- original row 58 data (columns B:AC)
$ws.Range("B60").Value = 7542499
$ws.Range("F60").Value = "Maccabi Petach Tikva"
$ws.Range("G60").Value = "Hapoel Beer Sheva"
$ws.Range("H60").Value = 1
$ws.Range("I60").Value = 4
$ws.Range("J60").Value = "A"
$ws.Range("K60").Value = 2.65
$ws.Range("L60").Value = 3.2
$ws.Range("M60").Value = 2.4
$ws.Range("N60").Value = 3.2
$ws.Range("O60").Value = 3.3
$ws.Range("P60").Value = 2.05
$ws.Range("Q60").Value = 0.25
$ws.Range("R60").Value = 2
$ws.Range("S60").Value = 1.85
$ws.Range("T60").Value = 2.25
$ws.Range("U60").Value = 1.85
$ws.Range("V60").Value = 2
$ws.Range("W60").Value = -1
$ws.Range("X60").Value = -1
$ws.Range("Y60").Value = 1.05
$ws.Range("Z60").Value = -1
$ws.Range("AA60").Value = 0.8500000000000001
$ws.Range("AB60").Value = 0.8500000000000001
$ws.Range("AC60").Value = -1

# Row 72 <- original row 73 data (columns B:AC)
$ws.Range("B72").Value = 7542639
$ws.Range("F72").Value = "Maccabi Bnei Raina"
$ws.Range("G72").Value = "Hapoel Jerusalem FC"
$ws.Range("H72").Value = 1
$ws.Range("I72").Value = 1
$ws.Range("J72").Value = "D"
$ws.Range("K72").Value = 2.5
$ws.Range("L72").Value = 3
$ws.Range("M72").Value = 2.75
$ws.Range("N72").Value = 2.7
$ws.Range("O72").Value = 2.8
$ws.Range("P72").Value = 2.75
$ws.Range("Q72").Value = 0
$ws.Range("R72").Value = 1.925
$ws.Range("S72").Value = 1.925
$ws.Range("T72").Value = 2
$ws.Range("U72").Value = 2.1
$ws.Range("V72").Value = 1.775
$ws.Range("W72").Value = -1
$ws.Range("X72").Value = 1.8
$ws.Range("Y72").Value = -1
$ws.Range("Z72").Value = 0
$ws.Range("AA72").Value = 0
$ws.Range("AB72").Value = 0
$ws.Range("AC72").Value = 0

# Row 73 <- original row 74 data (columns B:AC)
$ws.Range("B73").Value = 7542640
$ws.Range("F73").Value = "MS Ashdod"
$ws.Range("G73").Value = "Hapoel Bnei Sakhnin"
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 1
$ws.Range("J73").Value = "A"
$ws.Range("K73").Value = 2.05
$ws.Range("L73").Value = 3.2
$ws.Range("M73").Value = 3.5
$ws.Range("N73").Value = 2.15
$ws.Range("O73").Value = 3.1
$ws.Range("P73").Value = 3.2
$ws.Range("Q73").Value = -0.25
$ws.Range("R73").Value = 1.925
$ws.Range("S73").Value = 1.925
$ws.Range("T73").Value = 2.25
$ws.Range("U73").Value = 1.9
$ws.Range("V73").Value = 1.95
$ws.Range("W73").Value = -1
$ws.Range("X73").Value = -1
$ws.Range("Y73").Value = 2.2
$ws.Range("Z73").Value = -1
$ws.Range("AA73").Value = 0.925
$ws.Range("AB73").Value = -1
$ws.Range("AC73").Value = 0.95

# Row 74 <- original row 72 data (columns B:AC)
$ws.Range("B74").Value = 7542719
$ws.Range("F74").Value = "Hapoel Haifa"
$ws.Range("G74").Value = "Maccabi Netanya"
$ws.Range("H74").Value = 2
$ws.Range("I74").Value = 1
$ws.Range("J74").Value = "H"
$ws.Range("K74").Value = 2.6
$ws.Range("L74").Value = 3.1
$ws.Range("M74").Value = 2.6
$ws.Range("N74").Value = 2.9
$ws.Range("O74").Value = 3.2
$ws.Range("P74").Value = 2.3
$ws.Range("Q74").Value = 0.25
$ws.Range("R74").Value = 1.8
$ws.Range("S74").Value = 2.05
$ws.Range("T74").Value = 2.5
$ws.Range("U74").Value = 2
$ws.Range("V74").Value = 1.85
$ws.Range("W74").Value = 1.9
$ws.Range("X74").Value = -1
$ws.Range("Y74").Value = -1
$ws.Range("Z74").Value = 0.8
$ws.Range("AA74").Value = -1
$ws.Range("AB74").Value = 1
$ws.Range("AC74").Value = -1

# Row 86 <- original row 87 data (columns B:AC)
$ws.Range("B86").Value = 7542726
$ws.Range("F86").Value = "Hapoel Hadera"
$ws.Range("G86").Value = "Maccabi Netanya"
$ws.Range("H86").Value = 1
$ws.Range("I86").Value = 4
$ws.Range("J86").Value = "A"
$ws.Range("K86").Value = 3.3
$ws.Range("L86").Value = 3.5
$ws.Range("M86").Value = 2
$ws.Range("N86").Value = 4.333
$ws.Range("O86").Value = 3.6
$ws.Range("P86").Value = 1.7
$ws.Range("Q86").Value = 0.75
$ws.Range("R86").Value = 1.9
$ws.Range("S86").Value = 1.95
$ws.Range("T86").Value = 2.5
$ws.Range("U86").Value = 2
$ws.Range("V86").Value = 1.85
$ws.Range("W86").Value = -1
$ws.Range("X86").Value = -1
$ws.Range("Y86").Value = 0.7
$ws.Range("Z86").Value = -1
$ws.Range("AA86").Value = 0.95
$ws.Range("AB86").Value = 1
$ws.Range("AC86").Value = -1

# Row 87 <- original row 86 data (columns B:AC)
$ws.Range("B87").Value = 7542727
$ws.Range("F87").Value = "Maccabi Bnei Raina"
$ws.Range("G87").Value = "Hapoel Bnei Sakhnin"
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 1
$ws.Range("J87").Value = "A"
$ws.Range("K87").Value = 2.1
$ws.Range("L87").Value = 3.1
$ws.Range("M87").Value = 3.6
$ws.Range("N87").Value = 2.45
$ws.Range("O87").Value = 3
$ws.Range("P87").Value = 3
$ws.Range("Q87").Value = -0.25
$ws.Range("R87").Value = 2.075
$ws.Range("S87").Value = 1.725
$ws.Range("T87").Value = 2.25
$ws.Range("U87").Value = 2.05
$ws.Range("V87").Value = 1.8
$ws.Range("W87").Value = -1
$ws.Range("X87").Value = -1
$ws.Range("Y87").Value = 2
$ws.Range("Z87").Value = -1
$ws.Range("AA87").Value = 0.7250000000000001
$ws.Range("AB87").Value = -1
$ws.Range("AC87").Value = 0.8

# Row 109 <- original row 110 data (columns B:AC)
$ws.Range("B109").Value = 7542737
$ws.Range("F109").Value = "MS Ashdod"
$ws.Range("G109").Value = "Hapoel Haifa"
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 1
$ws.Range("J109").Value = "A"
$ws.Range("K109").Value = 3
$ws.Range("L109").Value = 3.2
$ws.Range("M109").Value = 2.45
$ws.Range("N109").Value = 3.2
$ws.Range("O109").Value = 3.25
$ws.Range("P109").Value = 2.3
$ws.Range("Q109").Value = 0.25
$ws.Range("R109").Value = 1.85
$ws.Range("S109").Value = 2
$ws.Range("T109").Value = 2.25
$ws.Range("U109").Value = 1.875
$ws.Range("V109").Value = 1.975
$ws.Range("W109").Value = -1
$ws.Range("X109").Value = -1
$ws.Range("Y109").Value = 1.3
$ws.Range("Z109").Value = -1
$ws.Range("AA109").Value = 1
$ws.Range("AB109").Value = -1
$ws.Range("AC109").Value = 0.9750000000000001

# Row 110 <- original row 109 data (columns B:AC)
$ws.Range("B110").Value = 7542736
$ws.Range("F110").Value = "Hapoel Jerusalem FC"
$ws.Range("G110").Value = "Hapoel Bnei Sakhnin"
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = "D"
$ws.Range("K110").Value = 2.2
$ws.Range("L110").Value = 3.4
$ws.Range("M110").Value = 3.2
$ws.Range("N110").Value = 2.375
$ws.Range("O110").Value = 3.1
$ws.Range("P110").Value = 3.1
$ws.Range("Q110").Value = -0.25
$ws.Range("R110").Value = 2.05
$ws.Range("S110").Value = 1.8
$ws.Range("T110").Value = 2
$ws.Range("U110").Value = 1.925
$ws.Range("V110").Value = 1.925
$ws.Range("W110").Value = -1
$ws.Range("X110").Value = 2.1
$ws.Range("Y110").Value = -1
$ws.Range("Z110").Value = -0.5
$ws.Range("AA110").Value = 0.4
$ws.Range("AB110").Value = -1
$ws.Range("AC110").Value = 0.925

# Row 144 <- original row 145 data (columns B:AC)
$ws.Range("B144").Value = 6799960
$ws.Range("F144").Value = "Maccabi Petach Tikva"
$ws.Range("G144").Value = "Maccabi Bnei Raina"
$ws.Range("H144").Value = 1
$ws.Range("I144").Value = 0
$ws.Range("J144").Value = "H"
$ws.Range("K144").Value = 2.625
$ws.Range("L144").Value = 3.25
$ws.Range("M144").Value = 2.5
$ws.Range("N144").Value = 2.8
$ws.Range("O144").Value = 3.25
$ws.Range("P144").Value = 2.375
$ws.Range("Q144").Value = 0.25
$ws.Range("R144").Value = 1.775
$ws.Range("S144").Value = 2.1
$ws.Range("T144").Value = 2.25
$ws.Range("U144").Value = 1.875
$ws.Range("V144").Value = 1.975
$ws.Range("W144").Value = 1.8
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = -1
$ws.Range("Z144").Value = 0.7749999999999999
$ws.Range("AA144").Value = -1
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = 0.9750000000000001

# Row 145 <- original row 144 data (columns B:AC)
$ws.Range("B145").Value = 6799962
$ws.Range("F145").Value = "MS Ashdod"
$ws.Range("G145").Value = "Hapoel Petah Tikva"
$ws.Range("H145").Value = 2
$ws.Range("I145").Value = 0
$ws.Range("J145").Value = "H"
$ws.Range("K145").Value = 2.2
$ws.Range("L145").Value = 3.1
$ws.Range("M145").Value = 3.2
$ws.Range("N145").Value = 2.2
$ws.Range("O145").Value = 3.1
$ws.Range("P145").Value = 3.2
$ws.Range("Q145").Value = -0.25
$ws.Range("R145").Value = 2
$ws.Range("S145").Value = 1.85
$ws.Range("T145").Value = 2.25
$ws.Range("U145").Value = 2
$ws.Range("V145").Value = 1.85
$ws.Range("W145").Value = 1.2
$ws.Range("X145").Value = -1
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = 1
$ws.Range("AA145").Value = -1
$ws.Range("AB145").Value = -0.5
$ws.Range("AC145").Value = 0.425

# Row 180 <- original row 181 data (columns B:AC)
$ws.Range("B180").Value = 6799995
$ws.Range("F180").Value = "Maccabi Petach Tikva"
$ws.Range("G180").Value = "Beitar Jerusalem"
$ws.Range("H180").Value = 0
$ws.Range("I180").Value = 3
$ws.Range("J180").Value = "A"
$ws.Range("K180").Value = 2.5
$ws.Range("L180").Value = 3.2
$ws.Range("M180").Value = 2.5
$ws.Range("N180").Value = 2.875
$ws.Range("O180").Value = 3.2
$ws.Range("P180").Value = 2.25
$ws.Range("Q180").Value = 0.25
$ws.Range("R180").Value = 1.825
$ws.Range("S180").Value = 2.025
$ws.Range("T180").Value = 2.5
$ws.Range("U180").Value = 1.975
$ws.Range("V180").Value = 1.875
$ws.Range("W180").Value = -1
$ws.Range("X180").Value = -1
$ws.Range("Y180").Value = 1.25
$ws.Range("Z180").Value = -1
$ws.Range("AA180").Value = 1.025
$ws.Range("AB180").Value = 0.9750000000000001
$ws.Range("AC180").Value = -1

# Row 181 <- original row 180 data (columns B:AC)
$ws.Range("B181").Value = 6799999
$ws.Range("F181").Value = "Hapoel Bnei Sakhnin"
$ws.Range("G181").Value = "Maccabi Bnei Raina"
$ws.Range("H181").Value = 0
$ws.Range("I181").Value = 0
$ws.Range("J181").Value = "D"
$ws.Range("K181").Value = 2.1
$ws.Range("L181").Value = 3
$ws.Range("M181").Value = 3.3
$ws.Range("N181").Value = 2.45
$ws.Range("O181").Value = 2.875
$ws.Range("P181").Value = 2.875
$ws.Range("Q181").Value = 0
$ws.Range("R181").Value = 1.8
$ws.Range("S181").Value = 2.05
$ws.Range("T181").Value = 2
$ws.Range("U181").Value = 2.025
$ws.Range("V181").Value = 1.825
$ws.Range("W181").Value = -1
$ws.Range("X181").Value = 1.875
$ws.Range("Y181").Value = -1
$ws.Range("Z181").Value = 0
$ws.Range("AA181").Value = 0
$ws.Range("AB181").Value = -1
$ws.Range("AC181").Value = 0.825
